$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45965
$ws.Range("B2").Value = 77.89
$ws.Range("C2").Value = 68.38
$ws.Range("D2").Value = 67.34999999999999
$ws.Range("E2").Value = 56.37
$ws.Range("F2").Value = 50.03
$ws.Range("G2").Value = 55.45
$ws.Range("H2").Value = 66.3
$ws.Range("I2").Value = 79
$ws.Range("J2").Value = 72.56
$ws.Range("K2").Value = 33.13
$ws.Range("L2").Value = 13.66
$ws.Range("M2").Value = 4.54
$ws.Range("N2").Value = 3.74
$ws.Range("O2").Value = 3.52
$ws.Range("P2").Value = 4.09
$ws.Range("Q2").Value = 12.3
$ws.Range("R2").Value = 23.23
$ws.Range("S2").Value = 57.37
$ws.Range("T2").Value = 77.47
$ws.Range("U2").Value = 83.56999999999999
$ws.Range("V2").Value = 82.64
$ws.Range("W2").Value = 78.12
$ws.Range("X2").Value = 72.34999999999999
$ws.Range("Y2").Value = 43.53
$ws.Range("Z2").Value = 49.44
$ws.Range("AB2").Value = 69.16
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 80.52
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 80.38
$ws.Range("AG2").Value = "9h-23h"
